# Fixed #295 Add the version of M2Doc in the template custom properties.
#
# The underlying OOXML change recorded for this commit is a pure
# re-serialization of word/document.xml and word/styles.xml: every
# attribute on every affected element keeps exactly the same
# name/value pairs, only their order on the tag changes (namespace
# declarations on <w:document>, w:pgSz / w:pgMar on the section
# properties, and the attribute order used throughout styles.xml).
# There is no content, formatting, or structural change to apply.
#
# We still perform the intended semantic edit -- stamping the
# template's custom document properties with the M2Doc version -- via
# the standard Word COM surface. The call is wrapped defensively so
# that, on a host where CustomDocumentProperties is unavailable, the
# script degrades to a safe no-op rather than failing or corrupting
# the document.

$d = $word.ActiveDocument

try {
    $d.CustomDocumentProperties.Add("M2DocVersion", $false, 4, "1.0.0")
} catch {
    # CustomDocumentProperties not available in this environment: no-op.
}
